{"js": "// Replace five paragraph texts (typo/name corrections + one appended\n// sentence) inside the Word document body, matching the target diff.\n\nconst replacements = [\n  {\n    find: \"Main- peaklass; psvm- meetod,  milles loome m\u00e4ngijate isendid ning k\u00e4ivitame m\u00e4ngu\",\n    text: \"Main- peaklass; psvm- meetod,  milles loome Mangijate isendid ning k\u00e4ivitame m\u00e4ngu\"\n  },\n  {\n    find: \"M\u00e4ngija- klass kuhu salvestame M\u00e4ngijate isendid; konstruktor ning getM\u00e4ngija - meetodid, mille abil loome M\u00e4ngija isendid ning v\u00e4ljastame neid hiljem ekraanile.\",\n    text: \"Mangija- klass kuhu salvestame Mangijate(m\u00e4ngijate) isendid; konstruktor ning getMangija - meetodid, mille abil loome Mangija isendid ning v\u00e4ljastame neid hiljem ekraanile.\"\n  },\n  {\n    find: \"M\u00e4ng- klass, kus kogu m\u00e4ng toimib; k\u00e4ik - meetod, milles m\u00e4ngija m\u00e4rk asetatakse m\u00e4ngulauale, kasNeliReas - meetod, mille abil kontrollime, kas neli samasugust m\u00e4rki on m\u00f5nes reas, veerus v\u00f5i diagonaalis, p\u00e4risM\u00e4ng - meetod, kus tehakse k\u00f5ik vajalikud tegevused, et m\u00e4ng toimiks.\",\n    text: \"Mang- klass, kus kogu m\u00e4ng toimib; kaik - meetod, milles m\u00e4ngija m\u00e4rk asetatakse m\u00e4ngulauale, kasNeliReas - meetod, mille abil kontrollime, kas neli samasugust m\u00e4rki on m\u00f5nes reas, veerus v\u00f5i diagonaalis, parisMang - meetod, kus tehakse k\u00f5ik vajalikud tegevused, et m\u00e4ng toimiks.\"\n  },\n  {\n    find: \"Sarah tegi meetodid t\u00e4idaTabel, k\u00e4ik, v\u00e4ljasta_tabel, ja p\u00e4risM\u00e4ng ja aitas ka \u00fclej\u00e4\u00e4nuga.\",\n    text: \"Sarah tegi meetodid taidaTabel, k\u00e4ik, valjastaTabel, ja parisMang ja aitas ka \u00fclej\u00e4\u00e4nuga.\"\n  },\n  {\n    find: \"Kindlasti oli k\u00f5ige raskem teha meetodit diagonaalide kontrollimise jaoks.Selleks otsisime abi googlest ning \u00fcritasime v\u00e4lja m\u00f5elda sobivat algoritmi, mis kontrolliks ilusti k\u00f5iki diagonaale. Samuti valmistas scanner klassi kasutamine esialgu segadust, kuid see kadus kiirelt.\",\n    text: \"Kindlasti oli k\u00f5ige raskem teha meetodit diagonaalide kontrollimise jaoks.Selleks otsisime abi googlest ning \u00fcritasime v\u00e4lja m\u00f5elda sobivat algoritmi, mis kontrolliks ilusti k\u00f5iki diagonaale. Samuti valmistas scanner klassi kasutamine esialgu segadust, kuid see kadus kiirelt. Samuti valmistas veidi raskusi \u00fchtse vormistuse ja stiili leidmine, kuna meil on v\u00e4lja kujunenud veidi erinev stiil.\"\n  }\n];\n\nfor (const { find, text } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(text, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply typo/name corrections (+ one appended sentence) via Find/Replace\n# on the whole document content, matching the target diff.\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{\n        Old = \"Main- peaklass; psvm- meetod,  milles loome m\u00e4ngijate isendid ning k\u00e4ivitame m\u00e4ngu\"\n        New = \"Main- peaklass; psvm- meetod,  milles loome Mangijate isendid ning k\u00e4ivitame m\u00e4ngu\"\n    },\n    @{\n        Old = \"M\u00e4ngija- klass kuhu salvestame M\u00e4ngijate isendid; konstruktor ning getM\u00e4ngija - meetodid, mille abil loome M\u00e4ngija isendid ning v\u00e4ljastame neid hiljem ekraanile.\"\n        New = \"Mangija- klass kuhu salvestame Mangijate(m\u00e4ngijate) isendid; konstruktor ning getMangija - meetodid, mille abil loome Mangija isendid ning v\u00e4ljastame neid hiljem ekraanile.\"\n    },\n    @{\n        Old = \"M\u00e4ng- klass, kus kogu m\u00e4ng toimib; k\u00e4ik - meetod, milles m\u00e4ngija m\u00e4rk asetatakse m\u00e4ngulauale, kasNeliReas - meetod, mille abil kontrollime, kas neli samasugust m\u00e4rki on m\u00f5nes reas, veerus v\u00f5i diagonaalis, p\u00e4risM\u00e4ng - meetod, kus tehakse k\u00f5ik vajalikud tegevused, et m\u00e4ng toimiks.\"\n        New = \"Mang- klass, kus kogu m\u00e4ng toimib; kaik - meetod, milles m\u00e4ngija m\u00e4rk asetatakse m\u00e4ngulauale, kasNeliReas - meetod, mille abil kontrollime, kas neli samasugust m\u00e4rki on m\u00f5nes reas, veerus v\u00f5i diagonaalis, parisMang - meetod, kus tehakse k\u00f5ik vajalikud tegevused, et m\u00e4ng toimiks.\"\n    },\n    @{\n        Old = \"Sarah tegi meetodid t\u00e4idaTabel, k\u00e4ik, v\u00e4ljasta_tabel, ja p\u00e4risM\u00e4ng ja aitas ka \u00fclej\u00e4\u00e4nuga.\"\n        New = \"Sarah tegi meetodid taidaTabel, k\u00e4ik, valjastaTabel, ja parisMang ja aitas ka \u00fclej\u00e4\u00e4nuga.\"\n    },\n    @{\n        Old = \"Kindlasti oli k\u00f5ige raskem teha meetodit diagonaalide kontrollimise jaoks.Selleks otsisime abi googlest ning \u00fcritasime v\u00e4lja m\u00f5elda sobivat algoritmi, mis kontrolliks ilusti k\u00f5iki diagonaale. Samuti valmistas scanner klassi kasutamine esialgu segadust, kuid see kadus kiirelt.\"\n        New = \"Kindlasti oli k\u00f5ige raskem teha meetodit diagonaalide kontrollimise jaoks.Selleks otsisime abi googlest ning \u00fcritasime v\u00e4lja m\u00f5elda sobivat algoritmi, mis kontrolliks ilusti k\u00f5iki diagonaale. Samuti valmistas scanner klassi kasutamine esialgu segadust, kuid see kadus kiirelt. Samuti valmistas veidi raskusi \u00fchtse vormistuse ja stiili leidmine, kuna meil on v\u00e4lja kujunenud veidi erinev stiil.\"\n    }\n)\n\nforeach ($rep in $replacements) {\n    $rng = $d.Content\n    $rng.Find.Execute($rep.Old, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, $rep.New, $wdReplaceAll)\n}\n"}
